$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 36.134107
$ws.Range("H2").Value = 108.402321
$ws.Range("I2").Value = 0.6652922211878534
$ws.Range("J2").Value = 0.6960141222607766
$ws.Range("M2").Value = 90.43008666666667
$ws.Range("N2").Value = 271.29026
$ws.Range("O2").Value = 0.863466363695901
$ws.Range("P2").Value = 0.8656179140344247
$ws.Range("Q2").Value = 3267.610427632607
$ws.Range("R2").Value = 29408.49384869346
$ws.Range("S2").Value = 0.5744574550242448
$ws.Range("T2").Value = 0.6024822926498745

$ws.Range("G3").Value = 36.134107
$ws.Range("H3").Value = 108.402321
$ws.Range("I3").Value = 0.6652922211878534
$ws.Range("J3").Value = 0.6960141222607766
$ws.Range("O3").Value = 0.000602137432244878
$ws.Range("P3").Value = 0.0006036378137891445
$ws.Range("Q3").Value = 2.278664966229333
$ws.Range("R3").Value = 20.507984696064
$ws.Range("S3").Value = 0.0004005973497585455
$ws.Range("T3").Value = 0.0004201404431278656

$ws.Range("G4").Value = 36.134107
$ws.Range("H4").Value = 108.402321
$ws.Range("I4").Value = 0.6652922211878534
$ws.Range("J4").Value = 0.6960141222607766
$ws.Range("M4").Value = 9.467965
$ws.Range("N4").Value = 28.403895
$ws.Range("O4").Value = 0.0904043069236993
$ws.Range("P4").Value = 0.09062957269587499
$ws.Range("Q4").Value = 342.116460382255
$ws.Range("R4").Value = 3079.048143440295
$ws.Range("S4").Value = 0.06014528215821634
$ws.Range("T4").Value = 0.06307946249078868

$ws.Range("G5").Value = 36.134107
$ws.Range("H5").Value = 108.402321
$ws.Range("I5").Value = 0.6652922211878534
$ws.Range("J5").Value = 0.6960141222607766
$ws.Range("M5").Value = 0.7809335
$ws.Range("N5").Value = 1.561867
$ws.Range("O5").Value = 0.007456697592460336
$ws.Range("P5").Value = 0.004983518592002547
$ws.Range("Q5").Value = 28.2183346488845
$ws.Range("R5").Value = 169.310007893307
$ws.Range("S5").Value = 0.004960882904014055
$ws.Range("T5").Value = 0.003468599318582914

$ws.Range("G6").Value = 36.134107
$ws.Range("H6").Value = 108.402321
$ws.Range("I6").Value = 0.6652922211878534
$ws.Range("J6").Value = 0.6960141222607766
$ws.Range("M6").Value = 3.98709
$ws.Range("N6").Value = 11.96127
$ws.Range("O6").Value = 0.03807049435569441
$ws.Range("P6").Value = 0.03816535686390858
$ws.Range("Q6").Value = 144.06993667863
$ws.Range("R6").Value = 1296.62943010767
$ws.Range("S6").Value = 0.02532800375161957
$ws.Range("T6").Value = 0.02656362735840264

$ws.Range("I7").Value = 0.1753444440667649
$ws.Range("J7").Value = 0.1834415095257413
$ws.Range("M7").Value = 90.43008666666667
$ws.Range("N7").Value = 271.29026
$ws.Range("O7").Value = 0.863466363695901
$ws.Range("P7").Value = 0.8656179140344247
$ws.Range("Q7").Value = 861.2115332372446
$ws.Range("R7").Value = 7750.903799135201
$ws.Range("S7").Value = 0.1514040295126088
$ws.Range("T7").Value = 0.1587902568229982

$ws.Range("I8").Value = 0.1753444440667649
$ws.Range("J8").Value = 0.1834415095257413
$ws.Range("O8").Value = 0.000602137432244878
$ws.Range("P8").Value = 0.0006036378137891445
$ws.Range("S8").Value = 0.0001055814533087674
$ws.Range("T8").Value = 0.000110732231768299

$ws.Range("I9").Value = 0.1753444440667649
$ws.Range("J9").Value = 0.1834415095257413
$ws.Range("M9").Value = 9.467965
$ws.Range("N9").Value = 28.403895
$ws.Range("O9").Value = 0.0904043069236993
$ws.Range("P9").Value = 0.09062957269587499
$ws.Range("Q9").Value = 90.16822779726667
$ws.Range("R9").Value = 811.5140501754
$ws.Range("S9").Value = 0.01585189293877724
$ws.Range("T9").Value = 0.01662522562300422

$ws.Range("I10").Value = 0.1753444440667649
$ws.Range("J10").Value = 0.1834415095257413
$ws.Range("M10").Value = 0.7809335
$ws.Range("N10").Value = 1.561867
$ws.Range("O10").Value = 0.007456697592460336
$ws.Range("P10").Value = 0.004983518592002547
$ws.Range("Q10").Value = 7.437225393473334
$ws.Range("R10").Value = 44.62335236084
$ws.Range("S10").Value = 0.001307490493923942
$ws.Range("T10").Value = 0.0009141841732665442

$ws.Range("I11").Value = 0.1753444440667649
$ws.Range("J11").Value = 0.1834415095257413
$ws.Range("M11").Value = 3.98709
$ws.Range("N11").Value = 11.96127
$ws.Range("O11").Value = 0.03807049435569441
$ws.Range("P11").Value = 0.03816535686390858
$ws.Range("Q11").Value = 37.97107819560001
$ws.Range("R11").Value = 341.7397037604001
$ws.Range("S11").Value = 0.006675449668146147
$ws.Range("T11").Value = 0.007001110674704002

$ws.Range("G12").Value = 0.5104573333333333
$ws.Range("H12").Value = 1.531372
$ws.Range("I12").Value = 0.00939841388954103
$ws.Range("J12").Value = 0.009832414367167749
$ws.Range("M12").Value = 90.43008666666667
$ws.Range("N12").Value = 271.29026
$ws.Range("O12").Value = 0.863466363695901
$ws.Range("P12").Value = 0.8656179140344247
$ws.Range("Q12").Value = 46.16070089296889
$ws.Range("R12").Value = 415.44630803672
$ws.Range("S12").Value = 0.008115214265711043
$ws.Range("T12").Value = 0.008511114014429855

$ws.Range("G13").Value = 0.5104573333333333
$ws.Range("H13").Value = 1.531372
$ws.Range("I13").Value = 0.00939841388954103
$ws.Range("J13").Value = 0.009832414367167749
$ws.Range("O13").Value = 0.000602137432244878
$ws.Range("P13").Value = 0.0006036378137891445
$ws.Range("Q13").Value = 0.03219012004977777
$ws.Range("R13").Value = 0.289711080448
$ws.Range("S13").Value = 0.000005659136806622833
$ws.Range("T13").Value = 0.000005935217112866115

$ws.Range("G14").Value = 0.5104573333333333
$ws.Range("H14").Value = 1.531372
$ws.Range("I14").Value = 0.00939841388954103
$ws.Range("J14").Value = 0.009832414367167749
$ws.Range("M14").Value = 9.467965
$ws.Range("N14").Value = 28.403895
$ws.Range("O14").Value = 0.0904043069236993
$ws.Range("P14").Value = 0.09062957269587499
$ws.Range("Q14").Value = 4.832992165993333
$ws.Range("R14").Value = 43.49692949393999
$ws.Range("S14").Value = 0.0008496570938660259
$ws.Range("T14").Value = 0.0008911075126651952

$ws.Range("G15").Value = 0.5104573333333333
$ws.Range("H15").Value = 1.531372
$ws.Range("I15").Value = 0.00939841388954103
$ws.Range("J15").Value = 0.009832414367167749
$ws.Range("M15").Value = 0.7809335
$ws.Range("N15").Value = 1.561867
$ws.Range("O15").Value = 0.007456697592460336
$ws.Range("P15").Value = 0.004983518592002547
$ws.Range("Q15").Value = 0.3986332319206666
$ws.Range("R15").Value = 2.391799391524
$ws.Range("S15").Value = 0.00007008113022308639
$ws.Range("T15").Value = 0.00004900001980305344

$ws.Range("G16").Value = 0.5104573333333333
$ws.Range("H16").Value = 1.531372
$ws.Range("I16").Value = 0.00939841388954103
$ws.Range("J16").Value = 0.009832414367167749
$ws.Range("M16").Value = 3.98709
$ws.Range("N16").Value = 11.96127
$ws.Range("O16").Value = 0.03807049435569441
$ws.Range("P16").Value = 0.03816535686390858
$ws.Range("Q16").Value = 2.03523932916
$ws.Range("R16").Value = 18.31715396244
$ws.Range("S16").Value = 0.0003578022629342518
$ws.Range("T16").Value = 0.000375257603156779

$ws.Range("G17").Value = 7.192107500000001
$ws.Range("H17").Value = 14.384215
$ws.Range("I17").Value = 0.1324193005156269
$ws.Range("J17").Value = 0.09235611087732429
$ws.Range("M17").Value = 90.43008666666667
$ws.Range("N17").Value = 271.29026
$ws.Range("O17").Value = 0.863466363695901
$ws.Range("P17").Value = 0.8656179140344247
$ws.Range("Q17").Value = 650.3829045409834
$ws.Range("R17").Value = 3902.2974272459
$ws.Range("S17").Value = 0.1143396118993831
$ws.Range("T17").Value = 0.07994510404596149

$ws.Range("G18").Value = 7.192107500000001
$ws.Range("H18").Value = 14.384215
$ws.Range("I18").Value = 0.1324193005156269
$ws.Range("J18").Value = 0.09235611087732429
$ws.Range("O18").Value = 0.000602137432244878
$ws.Range("P18").Value = 0.0006036378137891445
$ws.Range("Q18").Value = 0.4535438884266667
$ws.Range("R18").Value = 2.72126333056
$ws.Range("S18").Value = 0.00007973461759214246
$ws.Range("T18").Value = 0.00005574964086005586

$ws.Range("G19").Value = 7.192107500000001
$ws.Range("H19").Value = 14.384215
$ws.Range("I19").Value = 0.1324193005156269
$ws.Range("J19").Value = 0.09235611087732429
$ws.Range("M19").Value = 9.467965
$ws.Range("N19").Value = 28.403895
$ws.Range("O19").Value = 0.0904043069236993
$ws.Range("P19").Value = 0.09062957269587499
$ws.Range("Q19").Value = 68.09462208623751
$ws.Range("R19").Value = 408.567732517425
$ws.Range("S19").Value = 0.01197127508643631
$ws.Range("T19").Value = 0.008370194864664752

$ws.Range("G20").Value = 7.192107500000001
$ws.Range("H20").Value = 14.384215
$ws.Range("I20").Value = 0.1324193005156269
$ws.Range("J20").Value = 0.09235611087732429
$ws.Range("M20").Value = 0.7809335
$ws.Range("N20").Value = 1.561867
$ws.Range("O20").Value = 0.007456697592460336
$ws.Range("P20").Value = 0.004983518592002547
$ws.Range("Q20").Value = 5.61655768235125
$ws.Range("R20").Value = 22.466230729405
$ws.Range("S20").Value = 0.0009874106793501573
$ws.Range("T20").Value = 0.0004602583956421943

$ws.Range("G21").Value = 7.192107500000001
$ws.Range("H21").Value = 14.384215
$ws.Range("I21").Value = 0.1324193005156269
$ws.Range("J21").Value = 0.09235611087732429
$ws.Range("M21").Value = 3.98709
$ws.Range("N21").Value = 11.96127
$ws.Range("O21").Value = 0.03807049435569441
$ws.Range("P21").Value = 0.03816535686390858
$ws.Range("Q21").Value = 28.675579892175
$ws.Range("R21").Value = 172.05347935305
$ws.Range("S21").Value = 0.005041268232865178
$ws.Range("T21").Value = 0.003524803930195791

$ws.Range("E22").Value = 3
$ws.Range("F22").Value = 1
$ws.Range("G22").Value = 0.9529576666666667
$ws.Range("H22").Value = 2.858873
$ws.Range("I22").Value = 0.01754562034021376
$ws.Range("J22").Value = 0.01835584296898988
$ws.Range("M22").Value = 90.43008666666667
$ws.Range("N22").Value = 271.29026
$ws.Range("O22").Value = 0.863466363695901
$ws.Range("P22").Value = 0.8656179140344247
$ws.Range("Q22").Value = 86.17604438633111
$ws.Range("R22").Value = 775.58439947698
$ws.Range("S22").Value = 0.01515005299395322
$ws.Range("T22").Value = 0.01588914650116048

$ws.Range("E23").Value = 3
$ws.Range("F23").Value = 1
$ws.Range("G23").Value = 0.9529576666666667
$ws.Range("H23").Value = 2.858873
$ws.Range("I23").Value = 0.01754562034021376
$ws.Range("J23").Value = 0.01835584296898988
$ws.Range("O23").Value = 0.000602137432244878
$ws.Range("P23").Value = 0.0006036378137891445
$ws.Range("Q23").Value = 0.06009478107022222
$ws.Range("R23").Value = 0.540853029632
$ws.Range("S23").Value = 0.00001056487477879982
$ws.Range("T23").Value = 0.00001108028092005789

$ws.Range("E24").Value = 3
$ws.Range("F24").Value = 1
$ws.Range("G24").Value = 0.9529576666666667
$ws.Range("H24").Value = 2.858873
$ws.Range("I24").Value = 0.01754562034021376
$ws.Range("J24").Value = 0.01835584296898988
$ws.Range("M24").Value = 9.467965
$ws.Range("N24").Value = 28.403895
$ws.Range("O24").Value = 0.0904043069236993
$ws.Range("P24").Value = 0.09062957269587499
$ws.Range("Q24").Value = 9.022569834481667
$ws.Range("R24").Value = 81.20312851033499
$ws.Range("S24").Value = 0.001586199646403386
$ws.Range("T24").Value = 0.001663582204752134

$ws.Range("E25").Value = 3
$ws.Range("F25").Value = 1
$ws.Range("G25").Value = 0.9529576666666667
$ws.Range("H25").Value = 2.858873
$ws.Range("I25").Value = 0.01754562034021376
$ws.Range("J25").Value = 0.01835584296898988
$ws.Range("M25").Value = 0.7809335
$ws.Range("N25").Value = 1.561867
$ws.Range("O25").Value = 0.007456697592460336
$ws.Range("P25").Value = 0.004983518592002547
$ws.Range("Q25").Value = 0.7441965659818333
$ws.Range("R25").Value = 4.465179395891
$ws.Range("S25").Value = 0.0001308323849490951
$ws.Range("T25").Value = 0.00009147668470784027

$ws.Range("E26").Value = 3
$ws.Range("F26").Value = 1
$ws.Range("G26").Value = 0.9529576666666667
$ws.Range("H26").Value = 2.858873
$ws.Range("I26").Value = 0.01754562034021376
$ws.Range("J26").Value = 0.01835584296898988
$ws.Range("M26").Value = 3.98709
$ws.Range("N26").Value = 11.96127
$ws.Range("O26").Value = 0.03807049435569441
$ws.Range("P26").Value = 0.03816535686390858
$ws.Range("Q26").Value = 3.79952798319
$ws.Range("R26").Value = 34.19575184871
$ws.Range("S26").Value = 0.0006679704401292652
$ws.Range("T26").Value = 0.0007005572974493658
